$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters to indices: E=5, F=6, G=7, H=8
# List of (row, col, newValue) updates as described by the diff
$updates = @(
    @(5, 5, 67), @(5, 6, 36), @(5, 8, 36),
    @(6, 5, 22),
    @(10, 5, 203), @(10, 6, 87), @(10, 8, 87),
    @(11, 5, 147),
    @(12, 5, 221), @(12, 6, 111), @(12, 8, 111),
    @(13, 5, 76), @(13, 6, 35), @(13, 8, 35),
    @(14, 5, 67),
    @(15, 5, 93),
    @(16, 5, 86),
    @(17, 5, 42), @(17, 6, 19), @(17, 8, 19),
    @(20, 5, 53), @(20, 6, 20), @(20, 8, 20),
    @(21, 5, 72),
    @(22, 5, 87),
    @(23, 5, 98), @(23, 6, 42), @(23, 8, 42),
    @(25, 5, 96),
    @(26, 5, 59),
    @(27, 5, 142), @(27, 6, 69), @(27, 8, 69),
    @(30, 5, 101),
    @(31, 5, 42),
    @(32, 5, 101),
    @(33, 5, 131),
    @(34, 5, 99),
    @(35, 5, 65),
    @(38, 5, 46),
    @(39, 5, 112), @(39, 6, 40), @(39, 8, 40),
    @(40, 6, 56), @(40, 8, 56),
    @(42, 5, 152),
    @(43, 5, 47),
    @(44, 5, 140),
    @(45, 5, 57),
    @(46, 5, 122),
    @(47, 5, 204), @(47, 6, 85), @(47, 8, 85),
    @(48, 5, 102),
    @(49, 5, 114),
    @(50, 5, 97),
    @(51, 5, 99), @(51, 6, 35), @(51, 8, 35),
    @(52, 5, 6), @(52, 6, 4), @(52, 8, 4)
)

foreach ($u in $updates) {
    $row = $u[0]
    $col = $u[1]
    $val = $u[2]
    $ws.Cells.Item($row, $col).Value = $val
}
